# Correct a handful of mis-typed length-group range labels in column A.
# "100-110" -> "106-110", "150-160" -> "156-160",
# "160-170" -> "166-170", "170-180," -> "176-180,"
# (row order below matches the order the corrected labels were entered,
# which determines the order new shared-string entries are appended)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(22, 1).Value = "106-110"
$ws.Cells.Item(34, 1).Value = "166-170"
$ws.Cells.Item(36, 1).Value = "176-180,"
$ws.Cells.Item(32, 1).Value = "156-160"

# Update the view state left by the editor: scrolled so row 11 is at the
# top, with cell A33 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A33").Select() | Out-Null
